# Applies the "break out stock.yaml completed" update:
#  1. Appends three new rows (64-66) of stock data to the "day" sheet,
#     keeping the "bsecode" column (D) as text, matching the existing
#     convention used for prior rows on that sheet.
#  2. On the "month" sheet, converts the "bsecode" column (D) values for
#     rows 35-47 from text to numeric, keeping the same digits.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "day" sheet - append rows 64-66
# ---------------------------------------------------------------------
$dayWs = $wb.Worksheets.Item("day")

$newRows = @(
    @(1, "TVSMOTOR",   "Tvs Motor Company Limited", "532343", -0.43, 2354.6,  531888,  "day", "01/07/2024 11:34:38"),
    @(2, "OBEROIRLTY", "Oberoi Realty Limited",      "533273", 0.51,  1774.75, 764541,  "day", "01/07/2024 11:34:38"),
    @(3, "ZYDUSLIFE",  "Zydus Lifesciences Ltd",      "532321", -0.62, 1067.3,  1910647, "day", "01/07/2024 11:34:38")
)

$rowIndex = 64
foreach ($rowData in $newRows) {
    $dayWs.Cells.Item($rowIndex, 1).Value = $rowData[0]
    $dayWs.Cells.Item($rowIndex, 2).Value = $rowData[1]
    $dayWs.Cells.Item($rowIndex, 3).Value = $rowData[2]

    # Keep bsecode (column D) as text, same as the rest of this sheet's
    # historical rows - force text by pre-pending an apostrophe, then
    # strip the resulting quote-prefix style so no stray formatting is
    # introduced.
    $dayWs.Cells.Item($rowIndex, 4).Value = "'" + $rowData[3]
    $dayWs.Cells.Item($rowIndex, 4).Style = "Normal"

    $dayWs.Cells.Item($rowIndex, 5).Value = $rowData[4]
    $dayWs.Cells.Item($rowIndex, 6).Value = $rowData[5]
    $dayWs.Cells.Item($rowIndex, 7).Value = $rowData[6]
    $dayWs.Cells.Item($rowIndex, 8).Value = $rowData[7]
    $dayWs.Cells.Item($rowIndex, 9).Value = $rowData[8]

    $rowIndex = $rowIndex + 1
}

# ---------------------------------------------------------------------
# 2) "month" sheet - convert bsecode (column D) rows 35-47 to numbers
# ---------------------------------------------------------------------
$monthWs = $wb.Worksheets.Item("month")

for ($r = 35; $r -le 47; $r++) {
    $cell = $monthWs.Cells.Item($r, 4)
    $textValue = $cell.Value()
    $cell.Value = [double]$textValue
}

Write-Output "done"
